# Old Latex for Free-Stream Velocity
# Adds a new "Sheet2" (placed after "Sheet1") containing the free-stream
# velocity calibration calculations, and updates the selection on Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update the selection on the existing sheet (Sheet1)
$null = $ws1.Range("F8").Select()

# Insert a new worksheet right after Sheet1 and name it Sheet2
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Row 2: first value is literal, remaining cells reference the row below
$ws2.Range("A2").Value = 459
$ws2.Range("B2").Formula = "=A3"
$ws2.Range("C2:P2").Formula = "=B3"

# Row 3: literal calibration values
$row3 = @(545, 667, 766, 865, 967, 1069, 1166, 1257, 1346, 1442, 1554, 1630, 1746, 1852, 1962, 2048)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws2.Cells.Item(3, $i + 1).Value = $row3[$i]
}

# Row 4: difference between row 3 and row 2
$ws2.Range("A4").Formula = "=A3-A2"
$ws2.Range("B4:P4").Formula = "=B3-B2"

# Row 6: average, "+/-" label (typed with a leading apostrophe so Excel
# stores it as text with a quote prefix), and standard deviation
$ws2.Range("A6").Formula = "=AVERAGE(A4:P4)"
$ws2.Range("B6").Value = "'+/-"
$ws2.Range("C6").Formula = "=STDEV(A4:P4)"

# Row 7: relative uncertainty
$ws2.Range("C7").Formula = "=2*C6/A6"

# Row 9-10: free-stream velocity calibration
$ws2.Range("A9").Formula = "=9.68/A6"
$ws2.Range("A10").Formula = "=A9*206"

# Select E9 and make Sheet2 the active/visible tab
$null = $ws2.Range("E9").Select()
$null = $ws2.Activate()
